$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Thursday" date column (D) for the week, mirroring
# column C's values where applicable.

# D2: header date, one day after C2 (which is 43054 / 2017-11-15)
$ws.Range("D2").Value = 43055
$ws.Range("D2").NumberFormat = $ws.Range("C2").NumberFormat

# D3, D5, D7 get a value of 3 (hours), matching style of existing cells
$ws.Range("D3").Value = 3
$ws.Range("D5").Value = 3
$ws.Range("D7").Value = 3

# Update the active selection to E9
$ws.Range("E9").Select()
